$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "Evolutionary strategy 2 - Elitarian selection evolution"
$ws.Range("E28").Value = "https://ropiens.tistory.com/138"

$ws.Range("D44").Value = "Video Coding for Machine (VCM) 논문 리뷰"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/84"

$ws.Range("D46").Value = "맹장염 (충수염)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/410"
